$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (formerly the "LC / Dropptaggsvamp / Hydnellum ferrugineum" record)
# becomes the "NT / Blå taggsvamp / Hydnellum caeruleum" record (swapped with row 5),
# while column B gets a fresh, independent value.
$ws.Range("A4").Value = 112127587
$ws.Range("B4").Value = 90808
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 4362
$ws.Range("F4").Value = "Blå taggsvamp"
$ws.Range("G4").Value = "Hydnellum caeruleum"
$ws.Range("H4").Value = "(Hornem.) P.Karst."
$ws.Range("P4").Value = "Svarvarmyran (Svarvarmyran), Ly lm"
$ws.Range("Q4").Value = 690447
$ws.Range("R4").Value = 7125629
$ws.Range("Z4").Value = "14:27"
$ws.Range("AB4").Value = "14:27"

# Row 5 (formerly the "NT / Blå taggsvamp / Hydnellum caeruleum" record)
# becomes the "LC / Dropptaggsvamp / Hydnellum ferrugineum" record (swapped with row 4),
# while column B gets a fresh, independent value.
$ws.Range("A5").Value = 112128524
$ws.Range("B5").Value = 90814
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 4364
$ws.Range("F5").Value = "Dropptaggsvamp"
$ws.Range("G5").Value = "Hydnellum ferrugineum"
$ws.Range("H5").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("P5").Value = "Godmyr (Godmyr), Ly lm"
$ws.Range("Q5").Value = 690281
$ws.Range("R5").Value = 7126404
$ws.Range("Z5").Value = "15:00"
$ws.Range("AB5").Value = "15:00"

# Row 6: only column B changes.
$ws.Range("B6").Value = 90806
